# Applies the "answers" update for the three-digit ÷ one-digit worksheet.
# Each old text is unique within the document, so a simple
# Find/Replace (wdReplaceAll, but it will only ever match once) is safe.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-01-25 Saturday"; new = "2025-01-26 Sunday"},
    @{old = "553÷2=276, 1"; new = "237÷8=29, 5"},
    @{old = "466÷3=155, 1"; new = "822÷4=205, 2"},
    @{old = "682÷3=227, 1"; new = "265÷7=37, 6"},
    @{old = "715÷4=178, 3"; new = "418÷5=83, 3"},
    @{old = "834÷6=139, 0"; new = "206÷3=68, 2"},
    @{old = "763÷2=381, 1"; new = "309÷6=51, 3"},
    @{old = "577÷7=82, 3"; new = "634÷2=317, 0"},
    @{old = "889÷4=222, 1"; new = "648÷6=108, 0"},
    @{old = "774÷3=258, 0"; new = "149÷9=16, 5"},
    @{old = "871÷6=145, 1"; new = "782÷3=260, 2"},
    @{old = "764÷7=109, 1"; new = "286÷2=143, 0"},
    @{old = "843÷5=168, 3"; new = "908÷8=113, 4"},
    @{old = "309÷3=103, 0"; new = "811÷7=115, 6"},
    @{old = "701÷9=77, 8"; new = "850÷5=170, 0"},
    @{old = "122÷5=24, 2"; new = "750÷9=83, 3"},
    @{old = "651÷5=130, 1"; new = "399÷3=133, 0"},
    @{old = "779÷5=155, 4"; new = "311÷6=51, 5"},
    @{old = "261÷6=43, 3"; new = "657÷2=328, 1"},
    @{old = "553÷6=92, 1"; new = "410÷9=45, 5"},
    @{old = "372÷9=41, 3"; new = "706÷9=78, 4"},
    @{old = "265÷2=132, 1"; new = "259÷2=129, 1"},
    @{old = "689÷8=86, 1"; new = "458÷6=76, 2"},
    @{old = "158÷4=39, 2"; new = "838÷9=93, 1"},
    @{old = "962÷3=320, 2"; new = "368÷5=73, 3"},
    @{old = "303÷6=50, 3"; new = "202÷4=50, 2"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
